# Insert a new weekly price record at row 204 (pushing the existing
# rows 204-243 down to 205-244) on the single worksheet of the workbook.
#
# Net effect matches the target diff: dimension grows from A1:R243 to
# A1:R244, a brand-new row 204 is inserted with the latest week's data,
# and every previously existing data row (204-243) simply shifts down by
# one row (205-244) with its values unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 204:243 down to 205:244, leaving a blank row 204 behind.
$ws.Rows("204:204").Insert()

# Populate the newly-inserted row 204 with this week's record.
$ws.Range("A204").Value2 = 4
$ws.Range("B204").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value2 = "Los Lagos"
$ws.Range("D204").Value2 = 44641
$ws.Range("E204").Value2 = 10
$ws.Range("F204").Value2 = 100112043
$ws.Range("G204").Value2 = "Pepino ensalada"
$ws.Range("H204").Value2 = "Sin especificar"
$ws.Range("I204").Value2 = "Primera"
$ws.Range("J204").Value2 = 150
$ws.Range("K204").Value2 = 22000
$ws.Range("L204").Value2 = 22000
$ws.Range("M204").Value2 = 22000
$ws.Range("N204").Value2 = "`$/caja 60 unidades"
$ws.Range("O204").Value2 = "Región de Arica y Parinacota"
$ws.Range("P204").Value2 = 367
$ws.Range("Q204").Value2 = 60
$ws.Range("R204").Value2 = "Hortaliza"
